$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 141, shifting rows 141:209 down to 142:210
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new data point
$ws.Cells.Item(141, 1).Value = 3
$ws.Cells.Item(141, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(141, 3).Value = "Coquimbo"
$ws.Cells.Item(141, 4).Value = 44460
$ws.Cells.Item(141, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(141, 5).Value = 5
$ws.Cells.Item(141, 6).Value = 100112032
$ws.Cells.Item(141, 7).Value = "Zapallo italiano"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 78
$ws.Cells.Item(141, 11).Value = 11000
$ws.Cells.Item(141, 12).Value = 12000
$ws.Cells.Item(141, 13).Value = 11487
$ws.Cells.Item(141, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(141, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(141, 16).Value = 164
$ws.Cells.Item(141, 17).Value = 70
$ws.Cells.Item(141, 18).Value = "Hortaliza"
